$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells, copying the style from an existing header cell
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122) # xlPasteFormats

# Update existing MSE / R2 / MAE values
$ws.Range("B2").Value = 0.1048784994476929
$ws.Range("C2").Value = 0.9980841181889433
$ws.Range("D2").Value = 0.2374218334876481

$ws.Range("B3").Value = 0.2529092640451192
$ws.Range("C3").Value = 0.9813141422307591
$ws.Range("D3").Value = 0.3935443691105401

# Populate the new Elapsed Time / CPU columns
$ws.Range("G2").Value = 0.2668650318499809
$ws.Range("H2").Value = 0.998

$ws.Range("G3").Value = 0.2668650318499809
$ws.Range("H3").Value = 0.998
